$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.959.28'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '1.879.43'
$ws.Range("E3").Value = '  +1.33%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9971'
$ws.Range("E4").Value = '  -0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.57'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9980'
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5157'
$ws.Range("E7").Value = '  +1.58%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3712'
$ws.Range("E8").Value = '  +1.79%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07181'
$ws.Range("E9").Value = '  +0.61%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8992'
$ws.Range("E10").Value = '  +1.36%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.75'
$ws.Range("E11").Value = '  +0.17%  '

# Row 12
$ws.Range("E12").Value = '  +0.52%  '

# Row 13
$ws.Range("D13").Value = '1.874.60'
$ws.Range("E13").Value = '  +1.07%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.90'
$ws.Range("E14").Value = '  +3.74%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.245'
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9971'
$ws.Range("E16").Value = '  -0.41%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008484'
$ws.Range("E17").Value = '  -0.52%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.24'
$ws.Range("E18").Value = '  +1.33%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9978'
$ws.Range("E19").Value = '  -0.23%  '

# Row 20
$ws.Range("D20").Value = '26.994.83'
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.027'
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("D22").Value = '2.118.89'
$ws.Range("E22").Value = '  +1.52%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.40'
$ws.Range("E23").Value = '  +1.46%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.438'
$ws.Range("E24").Value = '  +0.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.57'
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$ws.Range("E26").Value = '  -2.17%  '

# Row 27
$ws.Range("E27").Value = '  +1.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.115'
$ws.Range("E28").Value = '  +3.17%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.46'
$ws.Range("E29").Value = '  +1.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.912'
$ws.Range("E30").Value = '  +5.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.773'
$ws.Range("E31").Value = '  +3.15%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09184'
$ws.Range("E32").Value = '  -0.93%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05030'
$ws.Range("E33").Value = '  -1.54%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7553'
$ws.Range("E34").Value = '  +3.37%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +2.58%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.997'
$ws.Range("E36").Value = '  -2.32%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.278'
$ws.Range("E37").Value = '  +2.86%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01991'
$ws.Range("E38").Value = '  -0.69%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5584'
$ws.Range("E39").Value = '  +5.82%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.489'
$ws.Range("E40").Value = '  +1.47%  '

# Row 41
$ws.Range("E41").Value = '  +0.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.578'
$ws.Range("E42").Value = '  +1.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.75'
$ws.Range("E43").Value = '  -0.49%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.751'
$ws.Range("E44").Value = '  +4.35%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1502'
$ws.Range("E45").Value = '  +2.13%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4787'
$ws.Range("E46").Value = '  +3.29%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.14'
$ws.Range("E47").Value = '  +2.04%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9978'
$ws.Range("E48").Value = '  -0.24%  '

# Row 49
$ws.Range("E49").Value = '  +0.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.06'
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.34'
$ws.Range("E51").Value = '  +0.70%  '
